# Applies the crypto price/volume/coin updates described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '58.618.57'
$ws.Range('E2').Value = '  +1.70%  '
$ws.Range('D3').Value = '2.982.59'
$ws.Range('E3').Value = '  +2.82%  '
$ws.Range('E4').Value = '  +0.20%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '558.86'
$ws.Range('E5').Value = '  +1.77%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '136.22'
$ws.Range('E6').Value = '  +10.51%  '
$ws.Range('E7').Value = '  +0.11%  '
$ws.Range('E8').Value = '  +4.01%  '
$ws.Range('D9').Value = '2.973.36'
$ws.Range('E9').Value = '  +2.54%  '
$ws.Range('E10').Value = '  +4.33%  '
$ws.Range('B11').Value = 'Cardano'
$ws.Range('C11').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.453'
$ws.Range('E11').Value = '  +3.30%  '
$ws.Range('B12').Value = 'Toncoin'
$ws.Range('C12').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.80'
$ws.Range('E12').Value = '  +1.21%  '
$ws.Range('E13').Value = '  +6.32%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '33.41'
$ws.Range('E14').Value = '  +3.24%  '
$ws.Range('E15').Value = '  +2.59%  '
$ws.Range('D16').Value = '3.474.12'
$ws.Range('E16').Value = '  +3.01%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '6.93'
$ws.Range('E17').Value = '  +5.31%  '
$ws.Range('D18').Value = '2.975.59'
$ws.Range('E18').Value = '  +2.96%  '
$ws.Range('D19').Value = '58.679.38'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '422.06'
$ws.Range('E20').Value = '  +3.64%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.46'
$ws.Range('E21').Value = '  +4.17%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.708'
$ws.Range('E22').Value = '  +5.64%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.07'
$ws.Range('E23').Value = '  +3.30%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '13.37'
$ws.Range('E24').Value = '  +3.64%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '80.10'
$ws.Range('E25').Value = '  +3.85%  '
$ws.Range('E26').Value = '  +0.05%  '
$ws.Range('E27').Value = '  +0.20%  '
$ws.Range('E28').Value = '  +8.02%  '
$ws.Range('E29').Value = '  +2.35%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.73'
$ws.Range('E30').Value = '  +6.83%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '25.60'
$ws.Range('E31').Value = '  +3.63%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.03'
$ws.Range('E32').Value = '  -0.12%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0992'
$ws.Range('E33').Value = '  +0.85%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.987'
$ws.Range('E34').Value = '  +8.18%  '
$ws.Range('D35').Value = '0.0₃0749'
$ws.Range('E35').Value = '  +20.45%  '
$ws.Range('E36').Value = '  +5.93%  '
$ws.Range('E37').Value = '  +2.26%  '
$ws.Range('E38').Value = '  +0.47%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '8.68'
$ws.Range('E39').Value = '  +2.45%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.74'
$ws.Range('E40').Value = '  +13.07%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '398.28'
$ws.Range('E41').Value = '  +10.69%  '
$ws.Range('D42').Value = '2.728.47'
$ws.Range('E42').Value = '  +4.33%  '
$ws.Range('E43').Value = '  +0.81%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0346'
$ws.Range('E44').Value = '  +0.92%  '
$ws.Range('E45').Value = '  +0.06%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '124.79'
$ws.Range('E46').Value = '  +3.92%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.241'
$ws.Range('E47').Value = '  +5.22%  '
$ws.Range('E48').Value = '  +3.19%  '
$ws.Range('E49').Value = '  +1.53%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '23.18'
$ws.Range('E50').Value = '  +1.25%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '31.79'
$ws.Range('E51').Value = '  +18.81%  '
